$wb = $excel.ActiveWorkbook

# --- Features sheet: rename "Features" column to "Feature", renumber CFS_n -> F_n ---
$ws1 = $wb.Worksheets.Item("Features")
$ws1.Range("B1").Value = "Feature"
$ws1.Range("A2").Value = "F_1"
$ws1.Range("C2").Value = "F_1"
$ws1.Range("A3").Value = "F_2"
$ws1.Range("C3").Value = "F_2"
$ws1.Range("A4").Value = "F_3"
$ws1.Range("C4").Value = "F_3"
$ws1.Range("A5").Value = "F_4"
$ws1.Range("C5").Value = "F_4"
$ws1.Rows("3").AutoFit() | Out-Null
$ws1.Range("D7").Select() | Out-Null

# --- Add a new "Tasks" sheet after "Features" ---
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Tasks"

$ws2.Range("A1").Value = "ID"
$ws2.Range("B1").Value = "Task"
$ws2.Range("C1").Value = "TODO"
$ws2.Range("D1").Value = "In Progress"
$ws2.Range("E1").Value = "Done"

$ws2.Range("A2").Value = "T_1"
$ws2.Range("B2").Value = "Take handle as input from the user"
$ws2.Range("E2").Value = "T_1"

$ws2.Range("A3").Value = "T_2"
$ws2.Range("B3").Value = "Call a basic api and print json data to console"
$ws2.Range("C3").Value = "T_2"

# Match the "wrap text / vertical top" body style used on the Features sheet,
# applied only to the populated cells so no stray blank cells get written.
$usedCells = @("A1","B1","C1","D1","E1","A2","B2","E2","A3","B3","C3")
foreach ($addr in $usedCells) {
    $ws2.Range($addr).WrapText = $true
    $ws2.Range($addr).VerticalAlignment = -4160
}

$tbl = $ws2.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws2.Range("A1:E3"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table2"
$tbl.TableStyle = "TableStyleLight8"

$ws2.Range("C3").Select() | Out-Null
